$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, 34700, 2.962002992630005, 0.10720547182764599, 0.1, 2.9786839508981395),
  @(3, 34790, 3.005122900009155, 0.10776025167327916, 0.1, 2.9848922656451236),
  @(4, 34881, 2.984125137329102, 0.1100331477774277, 0.1, 2.9726510656250498),
  @(5, 34973, 3.032219171524048, 0.11192658965250599, 0.1, 2.9486774044114292),
  @(6, 35065, 2.814259052276611, 0.11406584423812978, 0.1, 2.824782827331436),
  @(7, 35156, 2.669142007827759, 0.11637799996832815, 0.1, 2.710641255558845),
  @(8, 35247, 2.651050090789795, 0.11680590486883849, 0.1, 2.6519235996663033),
  @(9, 35339, 2.636411905288696, 0.11781828691272139, 0.1, 2.5925757104730316),
  @(10, 35431, 2.433090209960938, 0.11969441391035332, 0.1, 2.470989371621478),
  @(11, 35521, 2.418379783630371, 0.11991983430698828, 0.1, 2.40053542704056),
  @(12, 35612, 2.222220420837402, 0.1207674300078815, 0.1, 2.2924376046589794),
  @(13, 35704, 2.270013809204102, 0.12108744447355385, 0.1, 2.2781512755708553),
  @(14, 35796, 2.256533861160278, 0.12361358571305242, 0.1, 2.26786902224913),
  @(15, 35886, 2.243213176727295, 0.12572551635613766, 0.1, 2.2704910175464206),
  @(16, 35977, 2.350176334381104, 0.12804269543067673, 0.10049906312962073, 2.3122880582558802),
  @(17, 36069, 2.453269243240356, 0.13215661586979127, 0.10320404848322934, 2.3175652081816165),
  @(18, 36161, 2.032520294189453, 0.13627286241642483, 0.10481543746996388, 2.1150687988159698),
  @(19, 36251, 1.963053822517395, 0.13855803313515427, 0.1062963872689699, 2.0390304173728895),
  @(20, 36342, 2.066593647003174, 0.14103811655207632, 0.10715154328019005, 2.0657308936339778),
  @(21, 36434, 1.881415724754333, 0.14516043747139928, 0.10846970376083774, 2.058239724517728),
  @(22, 36526, 2.447355270385742, 0.14986149911256358, 0.10750979954045362, 2.3560091311917253),
  @(23, 36617, 2.548131227493286, 0.15192403123359066, 0.1070238738318133, 2.4853536475877265),
  @(24, 36708, 2.530933618545532, 0.15228849138405404, 0.10500065599041414, 2.52874263347356),
  @(25, 36800, 2.574150085449219, 0.15460649685291178, 0.10438886243626072, 2.5722022107796603),
  @(26, 36892, 2.611109495162964, 0.15800243874393943, 0.10411768205251913, 2.609511327362066),
  @(27, 36982, 2.705683946609497, 0.16312149903252793, 0.10505834389166341, 2.6688669322020613),
  @(28, 37073, 2.633024215698242, 0.16877163995271835, 0.10488555755789769, 2.648807671147002),
  @(29, 37165, 2.782319068908691, 0.17697502246392205, 0.10567479674494937, 2.678161278136921),
  @(30, 37257, 2.436383247375488, 0.1858259878272013, 0.10509768831734462, 2.4581767086061337),
  @(31, 37347, 2.258062839508057, 0.19256917723737815, 0.1058754841180874, 2.288494400630591),
  @(32, 37438, 2.244787216186523, 0.200141254875395, 0.10568417699713087, 2.201182453917358),
  @(33, 37530, 1.963913083076477, 0.20986602617756867, 0.10599075307119873, 1.9685868055076896),
  @(34, 37622, 1.744187712669373, 0.21608688135411178, 0.10637556088762276, 1.744150957832967),
  @(35, 37712, 1.472136259078979, 0.222671934857118, 0.1079754277600149, 1.494251142573726),
  @(36, 37803, 1.254570722579956, 0.227652505819831, 0.10911534665506842, 1.2986433108788895),
  @(37, 37895, 1.093175888061523, 0.23060637324725147, 0.11040357649703096, 1.2124536876523386),
  @(38, 37987, 1.558441519737244, 0.23542864281172693, 0.11070656578285581, 1.5495713619988867),
  @(39, 38078, 1.865288138389587, 0.23750604800682606, 0.11129080268760799, 1.8294797929651914),
  @(40, 38169, 1.961798191070557, 0.2360717033952236, 0.11176978559839784, 1.9799414329171436),
  @(41, 38261, 2.265710115432739, 0.2356776432351506, 0.11226857714637786, 2.211421555181341),
  @(42, 38353, 2.352944374084473, 0.2332076635846343, 0.11298961226864773, 2.272885080429497),
  @(43, 38443, 2.034587860107422, 0.232875215432847, 0.11463303486518486, 2.0691264719360962),
  @(44, 38534, 1.924052119255066, 0.23262902793681178, 0.11518053913481141, 1.9900621440063604),
  @(45, 38626, 2.114802122116089, 0.23404035559987696, 0.11627524147754612, 2.110638056762385),
  @(46, 38718, 2.098948955535889, 0.23580397513562723, 0.11760709616128324, 2.189160942885992),
  @(47, 38808, 2.642067670822144, 0.24035926087224208, 0.1180898880116423, 2.5827993256560733),
  @(48, 38899, 2.930945873260498, 0.23887732792441968, 0.11853453174667306, 2.8097853073980437),
  @(49, 38991, 2.61341381072998, 0.2369725118034452, 0.11820416516525792, 2.6116044732929073),
  @(50, 39083, 2.50513768196106, 0.23321614322594678, 0.11841985563301147, 2.479410755071851),
  @(51, 39173, 2.181643724441528, 0.23162370431522922, 0.11796248692335114, 2.2395670462446695),
  @(52, 39264, 2.101353168487549, 0.22836402992624177, 0.11787471141257078, 2.1834200856662775),
  @(53, 39356, 2.435364246368408, 0.22632791930322288, 0.11850498806592791, 2.376452248193828),
  @(54, 39448, 2.388523817062378, 0.22435953239971784, 0.11848016786712007, 2.3787369150563116),
  @(55, 39539, 2.391724824905396, 0.22636421921743521, 0.11924335550598432, 2.369473754282062),
  @(56, 39630, 2.438619136810303, 0.22936673772867372, 0.12049745094337043, 2.3214049881285868),
  @(57, 39722, 1.762461423873901, 0.23565242353560567, 0.12018825266300562, 1.8696432079485839),
  @(58, 39814, 1.787616729736328, 0.23507680114041452, 0.11999915150606107, 1.792468480171229),
  @(59, 39904, 1.711974382400513, 0.2381805964759922, 0.12090651650043972, 1.7007166862334553),
  @(60, 39995, 1.479837775230408, 0.24038941194722008, 0.1210906899541691, 1.5480832142591388),
  @(61, 40087, 1.823668837547302, 0.24507160041839554, 0.12035668688129574, 1.6710374919543496),
  @(62, 40179, 1.159204602241516, 0.25049301636883614, 0.11927253306478344, 1.2112069060128325),
  @(63, 40269, 0.9501993060112, 0.2512760537960327, 0.11757560848267784, 0.9802271553053874),
  @(64, 40360, 0.8143872022628784, 0.24870726838091445, 0.11666100573924307, 0.8453529980540715),
  @(65, 40452, 0.6618974208831787, 0.24783832630005553, 0.11543039438276284, 0.781383373921595),
  @(66, 40544, 1.20978057384491, 0.2510340228108102, 0.11402201634218995, 1.2101280247924024),
  @(67, 40634, 1.583680391311646, 0.2482047138322048, 0.11233034275435816, 1.5878603417988137),
  @(68, 40725, 1.987721681594849, 0.2402919579272461, 0.11096486884967316, 1.95631949276662),
  @(69, 40817, 2.276660919189453, 0.22665637900058166, 0.11026383376945098, 2.203656426088922),
  @(70, 40909, 2.248338460922241, 0.21399534949336102, 0.10798646535455367, 2.217260962848638),
  @(71, 41000, 2.192283153533936, 0.20390085735989705, 0.10646403904590196, 2.158752910117135),
  @(72, 41091, 2.008243560791016, 0.19545265906084763, 0.10544846198363435, 2.0131638332013253),
  @(73, 41183, 1.899697661399841, 0.1875216763714092, 0.10486406128846032, 1.9151352285439203),
  @(74, 41275, 1.889026284217834, 0.18056298110073543, 0.10450189710976124, 1.8535747491662644),
  @(75, 41365, 1.62309741973877, 0.17569247757559864, 0.10472649585097157, 1.7051689830658534),
  @(76, 41456, 1.751937031745911, 0.17018145080975372, 0.10435274447358231, 1.7366733690066165),
  @(77, 41548, 1.740855097770691, 0.16639811724589162, 0.10369174511999837, 1.7332770074763886),
  @(78, 41640, 1.645658254623413, 0.1637665376546547, 0.10456646027140357, 1.710232145729191),
  @(79, 41730, 1.922862410545349, 0.16296809973369464, 0.10392318701079584, 1.8266236574417047),
  @(80, 41821, 1.740945339202881, 0.16152191944474648, 0.10385983976107974, 1.7488262638131538),
  @(81, 41913, 1.622418165206909, 0.1610712052980035, 0.10386038558798188, 1.6860751966871794),
  @(82, 42005, 1.74537980556488, 0.16057354444459943, 0.10346943009461457, 1.749512771135375),
  @(83, 42095, 1.777263283729553, 0.16185269442637518, 0.1043063167238543, 1.8039686315802148),
  @(84, 42186, 1.897095680236816, 0.16396990742428946, 0.10469871559631039, 1.908082258341556),
  @(85, 42278, 2.071506261825562, 0.16768346252052285, 0.10615108939099246, 2.047137465622313),
  @(86, 42370, 2.142422437667847, 0.1708632203538405, 0.10717768346027393, 2.132243362333834),
  @(87, 42461, 2.262210845947266, 0.1751093492487499, 0.10959190868129531, 2.2196542359196787),
  @(88, 42552, 2.271121025085449, 0.17831441191196587, 0.11182870106792081, 2.231371072911683),
  @(89, 42644, 2.197124242782593, 0.18443376241925624, 0.11411956929316752, 2.166830245887985),
  @(90, 42736, 2.04584789276123, 0.1918982013884673, 0.1168860115869725, 2.0171988449333593),
  @(91, 42826, 1.699123024940491, 0.19961693491809185, 0.11979006283215568, 1.7735771595182457),
  @(92, 42917, 1.595390319824219, 0.20513448387956507, 0.12293959291159597, 1.6915999201492298),
  @(93, 43009, 1.770164847373962, 0.2128183844267882, 0.12479360327146405, 1.8211322110867014),
  @(94, 43101, 2.122782707214355, 0.22080801327269217, 0.12873389313791692, 2.074090389422704),
  @(95, 43191, 2.245511293411255, 0.22590978204079334, 0.1314171058068996, 2.191298272819445),
  @(96, 43282, 2.197686433792114, 0.23353329723858443, 0.135317300701302, 2.1940215292352354),
  @(97, 43374, 2.24852728843689, 0.24381462234190054, 0.1398633142910305, 2.2105475642491816),
  @(98, 43466, 2.054852485656738, 0.2573295425620815, 0.14533704436769557, 2.102625106381228),
  @(99, 43556, 2.121228933334351, 0.2725406594633346, 0.14950123778195645, 2.1420697059064313),
  @(100, 43647, 2.344321250915527, 0.29046144691665055, 0.15691925334258683, 2.2827331722104423),
  @(101, 43739, 2.253693580627441, 0.3138677485463096, 0.1656732126178812, 2.21800595737239),
  @(102, 43831, 2.123221397399902, 0.34103556862879303, 0.1761593309190034, 2.0429419930242987),
  @(103, 43922, 1.176671504974365, 0.37515646208562103, 0.18614513383641706, 1.383742775276847),
  @(104, 44013, 1.706722378730774, 0.4029909336448343, 0.19440742149283033, 1.6785152516790707),
  @(105, 44105, 1.633168458938599, 0.4304701965053579, 0.20415380392229118, 1.7030636689146166),
  @(106, 44197, 1.660865187644958, 0.46740399035873204, 0.21248868328084417, 1.9547123518226794),
  @(107, 44287, 4.408313751220703, 0.508941633643567, 0.21880412159517007, 4.075705208473048),
  @(108, 44378, 4.022122859954834, 0.5295007453617545, 0.2213378778241542, 4.163922021309195),
  @(109, 44470, 5.522685050964355, 0.5550479545276172, 0.22288690608067296, 5.436205553613123),
  @(110, 44562, 6.452416896820068, 0.5599621586248305, 0.22327252772050782, 6.276032621554159),
  @(111, 44652, 5.884917736053467, 0.5659433207425455, 0.2248166416850833, 5.961449992634254),
  @(112, 44743, 6.64294958114624, 0.5604472116684359, 0.2260272994746485, 6.469774920130197),
  @(113, 44835, 5.703855991363525, 0.5615113374128287, 0.22433704039482716, 5.774368652635783),
  @(114, 44927, 5.602568626403809, 0.5499951075815052, 0.22278910492229648, 5.638531922715578)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 1).NumberFormat = "m/d/yy h:mm"
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}
